# HourRegistration.xlsx - "Implemented full scene transitions after winning/losing"
#
# Fills in previously-blank Activity entries for the "Web game" learning goal
# table (left block, columns C-G) and the matching right block (columns Q-U),
# replaces a couple of placeholder rows with the real logged entries, and
# removes now-superseded duplicate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Left table (Learning goal 1 / "Web game"), columns C:G ---------------

# Row 15: fill in the Activity that was previously left blank.
$ws.Range("G15").Value = "Working on physics/contact materials"

# Row 16: fill in the Activity that was previously left blank.
$ws.Range("G16").Value = "level, UI movement"

# Row 17: fill in the Activity that was previously left blank.
$ws.Range("G17").Value = "Added tweening lib.  Implemented door functionality, replaced HTML UI with THREE text"

# Row 18: fill in the Activity that was previously left blank.
$ws.Range("G18").Value = "Cleaning code, working on scene loading"

# Row 19: this used to be a placeholder ("13.00 - 17.00", 4 hrs, no activity).
# Replace it with the real logged entry.
$ws.Range("D19").Value = "10.00 - 16.00"
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = "Implemented scene loading"

# Row 20: this used to be a duplicate placeholder row ("13.00 - 17.00", 4 hrs,
# no activity, dated 44593). Replace it with the real logged entry (and its
# correct date, 44594).
$ws.Range("C20").Value = 44594
$ws.Range("D20").Value = "12.00 - 15.00"
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = "Writing reflection"

# Row 21: this was an incomplete placeholder (date + from-to only, no hours
# logged). It's now removed entirely.
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

# --- Right table (Learning goal 2), columns Q:U ----------------------------

# Row 20: used to be an unfinished placeholder ("9.00 - 12.00", 3 hrs, no
# activity, dated 44592). Replace with the real logged entry (correct date
# 44593).
$ws.Range("Q20").Value = 44593
$ws.Range("R20").Value = "12.00 - 14.00"
$ws.Range("T20").Value = 2
$ws.Range("U20").Value = "Full scene transitions winning/losing"

# Row 21: date shifts to 44594; activity for this entry is filled in.
$ws.Range("Q21").Value = 44594
$ws.Range("U21").Value = "Writing reflection"

# Row 22: this was a duplicate placeholder row, now removed entirely.
$ws.Range("Q22").ClearContents()
$ws.Range("R22").ClearContents()
$ws.Range("T22").ClearContents()

# --- Selection / view state -------------------------------------------------
$ws.Range("G16").Select()
